$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B (old B/C shift right to C/D).
$ws.Columns.Item(2).Insert()

# New header (B1) and new long stat-query cell (B2).
$ws.Cells.Item(1, 2).Value = "StatQuery"

$statQuery = @'
MATCH (s:study) WITH COLLECT(DISTINCT(s.clinical_study_designation)) AS all_studies MATCH (d:demographic) WITH COLLECT(DISTINCT(d.breed)) AS all_breeds, COLLECT(DISTINCT(d.sex)) AS all_sexes, all_studies MATCH (d:diagnosis) WITH COLLECT(DISTINCT(d.disease_term)) AS all_diseases, all_breeds, all_sexes, all_studies MATCH (p:program)<-[*]-(s:study)<-[*]-(c:case)<--(demo:demographic), (c)<--(diag:diagnosis) WHERE s.clinical_study_type IN ['Clinical Trial']  OPTIONAL MATCH (f:file)-[*]->(c), (samp:sample)-[*]->(c) WITH DISTINCT c AS c, p, s, demo, diag, f, samp RETURN count(DISTINCT(f)) as number_of_files , count(DISTINCT(samp)) as number_of_sample , count(DISTINCT(c.case_id)) as number_of_cases , count(DISTINCT(s.clinical_study_designation)) as number_of_study
'@

$ws.Cells.Item(2, 2).Value = $statQuery

# The new column wraps text like the neighbouring A2 cell.
$ws.Cells.Item(2, 2).WrapText = $true

# Column widths: new col B matches col A; C/D keep their original (shifted) widths.
$ws.Columns.Item(2).ColumnWidth = $ws.Columns.Item(1).ColumnWidth

# Update the view: active/top-left cell moves to C2.
$ws.Range("C2").Select()
